$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected Hydrogen demand (B3) and drop the stray Non-metallic minerals
# figure that used to sit in D3.
$ws.Range("B3").Value = 1625037.056697357
$ws.Range("D3").ClearContents()

# Corrected Chemicals-column figures for Methanol / Ammonia.
$ws.Range("C4").Value = 33.93168714694506
$ws.Range("C5").Value = 2424.883513345702

# Row 7 used to be labelled "Other"; it is now "Biogas" with an updated value.
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 411.770166210904

# A brand-new row 8 carries the "Other" label (with its own corrected value),
# styled the same way as the other row headers in column A.
$ws.Range("A8").Value = "Other"
$ws.Range("A7").Copy() | Out-Null
$ws.Range("A8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("D8").Value = 126.8499407842672
